$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column O: "寻址方式（可选）" (addressing method - optional), with
# sample values "静态" (static) / "动态" (dynamic) for the first two data rows.
$ws.Range("O1").Value = "寻址方式（可选）"
$ws.Range("O2").Value = "静态"
$ws.Range("O3").Value = "动态"

# Match the new column width introduced alongside the data (~18.875 chars).
$ws.Range("O1").ColumnWidth = 18.14

# Reflect the author's final selection on the newly added cell.
$ws.Range("O3").Select()
